# feat: add single and multi corrector
# Apply cell-level text corrections to the lab report sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - WBC
$ws.Range("B2").Value = "中性细胞数"
$ws.Range("D2").Value = "10^9/L"
$ws.Range("E2").Value = "3.5-9.5"

# Row 3 - RBC
$ws.Range("B3").Value = "中性细胞数"
$ws.Range("D3").Value = "10^12/L"
$ws.Range("E3").Value = "3.8-5.1"

# Row 4 - HGB
$ws.Range("C4").Value = "103.0"
$ws.Range("E4").Value = "115-150"

# Row 5 - HCT
$ws.Range("C5").Value = "34.5"
$ws.Range("E5").Value = "35-45"

# Row 6 - MCV
$ws.Range("C6").Value = "81.0"
$ws.Range("D6").Value = "fL"
$ws.Range("E6").Value = "82-100"

# Row 7 - MCH
$ws.Range("C7").Value = "24.2"
$ws.Range("D7").Value = "pg"
$ws.Range("E7").Value = "27-34"

# Row 8 - MCHO
$ws.Range("E8").Value = "316-354"

# Row 9 - PLT
$ws.Range("B9").Value = "血小板计数"
$ws.Range("D9").Value = "10^9/L"
$ws.Range("E9").Value = "125-350"

# Row 10 - RDW-CV
$ws.Range("B10").Value = "红细胞分布宽度"
$ws.Range("C10").Value = "18.0"
$ws.Range("E10").Value = "11.5-14.5"

# Row 11 - PDW
$ws.Range("B11").Value = "血小板分布宽度"
$ws.Range("C11").Value = "11.7"
$ws.Range("E11").Value = "9-17"

# Row 12 - MPV
$ws.Range("C12").Value = "10.5"
$ws.Range("D12").Value = "fL"
$ws.Range("E12").Value = "6-11.5"

# Row 13 - PCT
$ws.Range("E13").Value = "0.101-0.36"

# Row 14 - NEUT#
$ws.Range("D14").Value = "10^9/L"
$ws.Range("E14").Value = "1.8-6.3"

# Row 15 - LYMPH#
$ws.Range("D15").Value = "10^9/L"
$ws.Range("E15").Value = "1.1-3.2"

# Row 16 - MONO#
$ws.Range("D16").Value = "10^9/L"

# Row 17 - E#
$ws.Range("D17").Value = "10^9/L"
$ws.Range("E17").Value = "0.02-0.52"

# Row 18 - BASO#
$ws.Range("D18").Value = "10^9/L"
$ws.Range("E18").Value = "0-0.06"

# Row 19 - NEUT%
$ws.Range("B19").Value = "中性粒细胞百分率"
$ws.Range("E19").Value = "40-75"

# Row 20 - YMPH9 -> YMPH
$ws.Range("A20").Value = "YMPH"
$ws.Range("B20").Value = "淋巴细胞百分率"
$ws.Range("E20").Value = "20-50"

# Row 21 - MONO%
$ws.Range("B21").Value = "单核细胞百分率"
$ws.Range("E21").Value = "1-10"

# Row 22 - E0%
$ws.Range("B22").Value = "嗜酸性粒细胞百分率"
$ws.Range("E22").Value = "0.4-8"

# Row 23 - BAS0%
$ws.Range("B23").Value = "嗜碱性粒细胞百分率"
$ws.Range("E23").Value = "0-1"

# Row 24 - P-LCR
$ws.Range("B24").Value = "大型血小板比率"
$ws.Range("E24").Value = "17.5-42.3"

# Row 25 - DW-SD
$ws.Range("B25").Value = "红细胞分布宽度"
$ws.Range("D25").Value = "%"
$ws.Range("E25").Value = "38.2-49.2"
